# Doing Updates for Financials
# Restated historical financial figures on the JSAIY income/balance-sheet/cash-flow sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSAIY")

# row => @{ column = new value }
$updates = @{
    8 = @{ "D"=37115400; "E"=34204200; "F"=30659100; "G"=31010000; "H"=31236900; "I"=30394300; "J"=29078300 }
    9 = @{ "D"=34660700; "E"=32073000; "F"=28760000; "G"=29434400; "H"=29427800; "I"=28728700; "J"=27498800 }
    10 = @{ "D"=2454700; "E"=2131200; "F"=1899100; "G"=1575600; "H"=1809100; "I"=1665600; "J"=1579500 }
    14 = @{ "E"=-127800 }
    17 = @{ "D"=36439800; "E"=33366900; "F"=29737000; "G"=30904300; "H"=29920900; "I"=29243900; "J"=27938300 }
    18 = @{ "D"=675600; "E"=837400; "F"=922100; "G"=105600; "H"=1316000; "I"=1150400; "J"=1140000 }
    20 = @{ "E"=-32600; "F"=-50900; "G"=-39100; "H"=23500; "I"=23500; "J"=63900 }
    21 = @{ "D"=1639700; "E"=1633000; "F"=1641500; "G"=830200; "H"=2066200; "I"=1855800; "J"=1862000 }
    22 = @{ "D"=142200; "E"=148700; "F"=156500; "G"=160400; "H"=168300; "I"=167000; "J"=161700 }
    23 = @{ "D"=533500; "E"=656100; "F"=714800; "G"=-93900; "H"=1171300; "I"=1006900; "J"=1042100 }
    24 = @{ "D"=130400; "E"=164300; "F"=100400; "G"=122600; "H"=237400; "I"=221700; "J"=262200 }
    26 = @{ "D"=403000; "E"=491700; "F"=614300; "G"=-216500; "H"=933900; "I"=785200; "J"=780000 }
    27 = @{ "D"=379600; "E"=468200; "F"=598700; "G"=-216500; "H"=933900; "I"=785200; "J"=780000 }
    32 = @{ "E"=32600; "F"=50900; "G"=39100; "H"=-23500; "I"=-23500; "J"=-63900 }
    33 = @{ "D"=379600; "E"=468200; "F"=598700; "G"=-216500; "H"=933900; "I"=785200; "J"=780000 }
    35 = @{ "D"=379600; "E"=468200; "F"=598700; "G"=-216500; "H"=933900; "I"=785200; "J"=780000 }
    41 = @{ "D"=1365600; "E"=1098200; "F"=1113900; "G"=1676000; "H"=2076500; "I"=674300; "J"=963900 }
    42 = @{ "D"=1155600; "E"=444800; "F"=439600 }
    43 = @{ "D"=5148100; "E"=8300600; "F"=2733800; "G"=2570800; "H"=2158600; "I"=735600; "J"=523000 }
    44 = @{ "D"=2360800; "E"=4630300; "F"=1262600; "G"=1300400; "H"=1310800; "I"=2574700; "J"=2446900 }
    45 = @{ "D"=217800; "E"=326100; "F"=206100; "G"=219100; "H"=143500; "I"=198300; "J"=169600 }
    46 = @{ "D"=10248000; "E"=8232800; "F"=5755900; "G"=5766400; "H"=5689400; "I"=2479500; "J"=2650400 }
    47 = @{ "D"=4094200; "E"=6041600; "F"=3140800; "G"=2658200; "H"=2578600; "I"=990000; "J"=1020000 }
    48 = @{ "D"=12910100; "E"=26101900; "F"=12735300; "G"=12584000; "H"=12886600; "I"=18832900; "J"=24335800 }
    49 = @{ "D"=1398200; "E"=2015200; "F"=429100; "G"=423900; "H"=373000; "I"=423900; "J"=417400 }
    52 = @{ "D"=45700; "E"=39100; "F"=77000; "G"=137000; "H"=45700; "I"=78300; "J"=48300 }
    54 = @{ "D"=28696100; "E"=25822700; "F"=22138100; "G"=21569400; "H"=21573300; "I"=16558200; "J"=16095200 }
    57 = @{ "D"=3719900; "E"=13969200; "F"=2715600; "G"=2724700; "H"=2407800; "I"=6044200; "J"=3573800 }
    58 = @{ "D"=995200; "E"=606500; "F"=444800; "G"=430400; "H"=696500; "I"=242600; "J"=391300 }
    59 = @{ "D"=8721900; "E"=7314600; "F"=5604600; "G"=5874600; "H"=5719400; "I"=1359100; "J"=1412600 }
    60 = @{ "D"=13437000; "E"=11198800; "F"=8765000; "G"=9029700; "H"=8823700; "I"=4062900; "J"=4090300 }
    61 = @{ "D"=2333400; "E"=2659500; "F"=3092500; "G"=3370300; "H"=2934700; "I"=3413400; "J"=3413400 }
    62 = @{ "D"=3259500; "E"=3980800; "F"=1973400; "G"=1944700; "H"=1982600; "I"=1693000; "J"=1129500 }
    66 = @{ "D"=19029900; "E"=16859500; "F"=13836100; "G"=14344800; "H"=13743500; "I"=8945000; "J"=8633200 }
    72 = @{ "D"=6586800; "E"=6067700; "F"=5320300; "G"=4921200; "H"=5532900; "I"=4519400; "J"=4771200 }
    76 = @{ "D"=9666200; "E"=8963200; "F"=8301900; "G"=7224600; "H"=7829800; "I"=7613300; "J"=7462000 }
    81 = @{ "D"=379600; "E"=468200; "F"=598700; "G"=-216500; "H"=933900; "I"=785200; "J"=780000 }
    83 = @{ "D"=953500; "E"=819100; "F"=761700; "G"=755200; "H"=718700; "I"=674300; "J"=650900 }
    89 = @{ "D"=1780400; "E"=1503900; "F"=511300; "G"=1188200; "H"=1224700; "I"=1279500; "J"=1391700 }
    91 = @{ "D"=-731700; "E"=-826900; "F"=-842600; "G"=-1240400; "H"=-1194700; "I"=-1391700; "J"=-1600400 }
    94 = @{ "D"=-613000; "E"=-978200; "F"=-521700; "G"=-1173900; "H"=555600; "I"=-1124300; "J"=-1151700 }
    96 = @{ "D"=-306500; "E"=-330000; "F"=-310400; "G"=-430400; "H"=-417400; "I"=-401700; "J"=-371700 }
    100 = @{ "D"=-318300; "E"=-607800; "F"=-167000; "G"=-409600; "H"=-378200; "I"=-461700; "J"=71700 }
    102 = @{ "D"=849100; "E"=-82200; "F"=-177400; "G"=-395200; "H"=1402100; "I"=-306500; "J"=311700 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
